$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.166.34"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "2.477.86"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.18"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.20"
$ws.Range("E6").Value = "  +4.09%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "2.478.31"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.94"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "2.944.13"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.52"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "67.054.56"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "2.450.56"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.01"
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.22"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.99"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.69"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.23"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.80"
$ws.Range("E26").Value = "  +3.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.35"
$ws.Range("E27").Value = "  +5.17%  "
$ws.Range("D28").Value = "2.589.90"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "0.0₃0904"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "511.85"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.72"
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.77"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.53"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.17"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.70"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.329"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.84"
$ws.Range("E44").Value = "  +2.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").Value = "  +3.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.44"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.49"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.515"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").Value = "0.0₆0256"
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0734"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.57"
$ws.Range("E51").Value = "  -0.80%  "
